$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = "29.475.35"
    "E2" = "  +0.78%  "
    "D3" = "1.880.45"
    "E3" = "  +1.24%  "
    "D4" = "1.001"
    "E4" = "  +0.05%  "
    "D5" = "0.7164"
    "E5" = "  +1.41%  "
    "D6" = "242.25"
    "E6" = "  +1.78%  "
    "D7" = "1.001"
    "E7" = "  +0.01%  "
    "D8" = "0.07855"
    "E8" = "  -1.54%  "
    "D9" = "0.3121"
    "E9" = "  +3.39%  "
    "D10" = "25.12"
    "E10" = "  +7.06%  "
    "D11" = "0.08262"
    "E11" = "  +1.09%  "
    "D12" = "1.886.80"
    "E12" = "  +1.81%  "
    "D13" = "0.7322"
    "E13" = "  +3.87%  "
    "D14" = "5.289"
    "E14" = "  +1.92%  "
    "D15" = "91.17"
    "E15" = "  +1.70%  "
    "D16" = "29.559.11"
    "E16" = "  +0.99%  "
    "D17" = "5.944"
    "E17" = "  +2.51%  "
    "D18" = "247.77"
    "E18" = "  +4.00%  "
    "D19" = "0.000007905"
    "E19" = "  -0.37%  "
    "E20" = "  +0.75%  "
    "B21" = "Dai"
    "C21" = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
    "D21" = "0.9994"
    "E21" = "  +0.04%  "
    "B22" = "Chainlink"
    "C22" = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
    "D22" = "8.011"
    "E22" = "  +7.16%  "
    "E23" = "  +0.10%  "
    "D24" = "0.1593"
    "E24" = "  +11.28%  "
    "D25" = "164.20"
    "E25" = "  +0.83%  "
    "D26" = "9.042"
    "E26" = "  +1.98%  "
    "E27" = "  +1.34%  "
    "D28" = "1.363"
    "E28" = "  -4.60%  "
    "D29" = "1.496"
    "E29" = "  +1.39%  "
    "D30" = "4.374"
    "E30" = "  +0.15%  "
    "D31" = "4.139"
    "E31" = "  +2.98%  "
    "D32" = "0.05313"
    "E32" = "  +2.53%  "
    "D33" = "1.936"
    "E33" = "  +0.56%  "
    "E34" = "  +3.80%  "
    "D35" = "0.7240"
    "E35" = "  +1.54%  "
    "D36" = "2.679"
    "E36" = "  +1.06%  "
    "D37" = "0.01870"
    "E37" = "  +1.06%  "
    "D38" = "1.260.31"
    "E38" = "  +10.76%  "
    "D39" = "2.731"
    "E39" = "  +0.23%  "
    "D40" = "0.9113"
    "E40" = "  -2.65%  "
    "D41" = "74.38"
    "E41" = "  +5.79%  "
    "D42" = "6.107"
    "E42" = "  +2.76%  "
    "D43" = "1.001"
    "E43" = "  +0.02%  "
    "D44" = "103.88"
    "E44" = "  +1.16%  "
    "D45" = "2.031.85"
    "E45" = "  +1.63%  "
    "D46" = "0.5328"
    "D47" = "1.772"
    "E47" = "  +0.71%  "
    "D48" = "2.921"
    "E48" = "  +13.15%  "
    "E49" = "  +0.08%  "
    "E50" = "  +1.85%  "
    "D51" = "9.272"
    "E51" = "  +1.18%  "
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

Write-Host "Applied" $updates.Count "cell updates"